{"js": "// Update the division-practice answers in the table to the newly\n// generated values. Each old \"a\u00f7b=c, r\" string is replaced with its\n// corresponding new value using body.search/insertText (\"Replace\").\n\nconst replacements = [\n    [\"406\u00f74=101, 2\", \"423\u00f78=52, 7\"],\n    [\"410\u00f76=68, 2\", \"430\u00f75=86, 0\"],\n    [\"978\u00f78=122, 2\", \"526\u00f76=87, 4\"],\n    [\"848\u00f73=282, 2\", \"809\u00f76=134, 5\"],\n    [\"815\u00f78=101, 7\", \"162\u00f76=27, 0\"],\n    [\"331\u00f77=47, 2\", \"269\u00f76=44, 5\"],\n    [\"997\u00f73=332, 1\", \"728\u00f74=182, 0\"],\n    [\"195\u00f72=97, 1\", \"626\u00f74=156, 2\"],\n    [\"962\u00f75=192, 2\", \"944\u00f74=236, 0\"],\n    [\"508\u00f77=72, 4\", \"817\u00f72=408, 1\"],\n    [\"480\u00f77=68, 4\", \"925\u00f76=154, 1\"],\n    [\"831\u00f76=138, 3\", \"196\u00f74=49, 0\"],\n    [\"522\u00f73=174, 0\", \"381\u00f78=47, 5\"],\n    [\"272\u00f72=136, 0\", \"528\u00f73=176, 0\"],\n    [\"906\u00f72=453, 0\", \"336\u00f78=42, 0\"],\n    [\"648\u00f73=216, 0\", \"429\u00f76=71, 3\"],\n    [\"470\u00f75=94, 0\", \"941\u00f77=134, 3\"],\n    [\"872\u00f78=109, 0\", \"383\u00f77=54, 5\"],\n    [\"989\u00f76=164, 5\", \"384\u00f74=96, 0\"],\n    [\"828\u00f75=165, 3\", \"199\u00f75=39, 4\"],\n    [\"543\u00f75=108, 3\", \"763\u00f74=190, 3\"],\n    [\"178\u00f76=29, 4\", \"218\u00f74=54, 2\"],\n    [\"438\u00f76=73, 0\", \"165\u00f76=27, 3\"],\n    [\"725\u00f73=241, 2\", \"367\u00f72=183, 1\"],\n    [\"587\u00f78=73, 3\", \"316\u00f77=45, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n    const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    results.load(\"items\");\n    await context.sync();\n\n    for (let i = 0; i < results.items.length; i++) {\n        results.items[i].insertText(newText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n}\n", "ps1": "# Update the division-practice answers in the table to the newly\n# generated values. Each old \"a\u00f7b=c, r\" string is replaced with its\n# corresponding new value using Find/Replace against the whole document.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"406\u00f74=101, 2\"; New = \"423\u00f78=52, 7\" },\n    @{ Old = \"410\u00f76=68, 2\";  New = \"430\u00f75=86, 0\" },\n    @{ Old = \"978\u00f78=122, 2\"; New = \"526\u00f76=87, 4\" },\n    @{ Old = \"848\u00f73=282, 2\"; New = \"809\u00f76=134, 5\" },\n    @{ Old = \"815\u00f78=101, 7\"; New = \"162\u00f76=27, 0\" },\n    @{ Old = \"331\u00f77=47, 2\";  New = \"269\u00f76=44, 5\" },\n    @{ Old = \"997\u00f73=332, 1\"; New = \"728\u00f74=182, 0\" },\n    @{ Old = \"195\u00f72=97, 1\";  New = \"626\u00f74=156, 2\" },\n    @{ Old = \"962\u00f75=192, 2\"; New = \"944\u00f74=236, 0\" },\n    @{ Old = \"508\u00f77=72, 4\";  New = \"817\u00f72=408, 1\" },\n    @{ Old = \"480\u00f77=68, 4\";  New = \"925\u00f76=154, 1\" },\n    @{ Old = \"831\u00f76=138, 3\"; New = \"196\u00f74=49, 0\" },\n    @{ Old = \"522\u00f73=174, 0\"; New = \"381\u00f78=47, 5\" },\n    @{ Old = \"272\u00f72=136, 0\"; New = \"528\u00f73=176, 0\" },\n    @{ Old = \"906\u00f72=453, 0\"; New = \"336\u00f78=42, 0\" },\n    @{ Old = \"648\u00f73=216, 0\"; New = \"429\u00f76=71, 3\" },\n    @{ Old = \"470\u00f75=94, 0\";  New = \"941\u00f77=134, 3\" },\n    @{ Old = \"872\u00f78=109, 0\"; New = \"383\u00f77=54, 5\" },\n    @{ Old = \"989\u00f76=164, 5\"; New = \"384\u00f74=96, 0\" },\n    @{ Old = \"828\u00f75=165, 3\"; New = \"199\u00f75=39, 4\" },\n    @{ Old = \"543\u00f75=108, 3\"; New = \"763\u00f74=190, 3\" },\n    @{ Old = \"178\u00f76=29, 4\";  New = \"218\u00f74=54, 2\" },\n    @{ Old = \"438\u00f76=73, 0\";  New = \"165\u00f76=27, 3\" },\n    @{ Old = \"725\u00f73=241, 2\"; New = \"367\u00f72=183, 1\" },\n    @{ Old = \"587\u00f78=73, 3\";  New = \"316\u00f77=45, 1\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($pair.Old, $false, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n"}
